$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "De'Aaron Fox"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Sacramento Kings"
$ws.Range("A3").Value = "Justin Champagnie"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Washington Wizards"
$ws.Range("A4").Value = "Scottie Barnes"
$ws.Range("B4").Value = "PG,SG,SF,PF"
$ws.Range("C4").Value = "Toronto Raptors"
$ws.Range("A5").Value = "Luke Kennard"
$ws.Range("B5").Value = "SG"
$ws.Range("C5").Value = "Memphis Grizzlies"
$ws.Range("A6").Value = "Miles Bridges"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Charlotte Hornets"
$ws.Range("A7").Value = "Mikal Bridges"
$ws.Range("B7").Value = "SG,SF,PF"
$ws.Range("C7").Value = "New York Knicks"
$ws.Range("A8").Value = "DeMar DeRozan"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Sacramento Kings"
$ws.Range("A9").Value = "Evan Mobley"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Cleveland Cavaliers"
$ws.Range("A10").Value = "Onyeka Okongwu"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Atlanta Hawks"
$ws.Range("A11").Value = "Goga Bitadze"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Orlando Magic"
$ws.Range("A12").Value = "Brook Lopez"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Milwaukee Bucks"
$ws.Range("A13").Value = "Draymond Green"
$ws.Range("B13").Value = "PF,C"
$ws.Range("C13").Value = "Golden State Warriors"
$ws.Range("A14").Value = "Tyler Herro"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Miami Heat"
$ws.Range("A15").Value = "Josh Giddey"
$ws.Range("B15").Value = "PG,SG,SF"
$ws.Range("C15").Value = "Chicago Bulls"
$ws.Range("A16").Value = "Nikola Vucevic"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Chicago Bulls"
$ws.Range("A17").Value = "Santi Aldama"
$ws.Range("B17").Value = "PF,C"
$ws.Range("C17").Value = "Memphis Grizzlies"
$ws.Range("A18").Value = "Luka Doncic"
$ws.Range("B18").Value = "PG,SG"
$ws.Range("C18").Value = "Dallas Mavericks"
$ws.Range("A19").Value = "Ja Morant"
$ws.Range("B19").Value = "PG"
$ws.Range("C19").Value = "Memphis Grizzlies"
